# "Su dung regular expression" - use a wildcard (regular-expression-like)
# Find & Replace to bump the team number in "Nhom thuc hien: Nhom 8" to 9.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "(Nhóm thực hiện: Nhóm )[0-9]{1,}",  # FindText (wildcard pattern, group 1 = prefix)
    $false,                               # MatchCase
    $false,                               # MatchWholeWord
    $true,                                # MatchWildcards
    $false,                               # MatchSoundsLike
    $false,                               # MatchAllWordForms
    $true,                                # Forward
    1,                                    # Wrap (wdFindContinue)
    $false,                               # Format
    "\1" + "9",                           # ReplaceWith (keep prefix, swap in 9)
    2                                     # Replace (wdReplaceAll)
) | Out-Null

# Also flips the (quasi-cosmetic) "Normal Table" style to a QuickStyle, as
# seen in styles.xml for this revision.
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true
